$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shorten / update header names (row 1)
$ws.Range("B1").Value = "Niketan (TL)"
$ws.Range("F1").Value = "Vishwajeet J"
$ws.Range("I1").Value = "Ashutosh W"
$ws.Range("J1").Value = "Uday W"
$ws.Range("K1").Value = "RaviKumar S"

# 2. Highlight the Team Lead header cell in yellow
$ws.Range("B1").Interior.Color = 65535

# 3. Normalize the font used for the TL's attendance cell on row 2
$ws.Range("B2").Font.Name = "Calibri"

# 4. Clear the stray ABSENT mark that used to sit under the date column
$ws.Range("B3").ClearContents()

# 5. Resize columns to fit the shorter header text
$ws.Columns.Item(2).ColumnWidth = 13.75
$ws.Columns.Item(6).ColumnWidth = 13.125
$ws.Columns.Item(7).ColumnWidth = 13.9375
$ws.Columns.Item(9).ColumnWidth = 12.125
$ws.Columns.Item(10).ColumnWidth = 11.75
$ws.Columns.Item(11).ColumnWidth = 14.75

# 6. Reset the view so column A is visible again and B1 is the active cell
$ws.Range("B1").Select() | Out-Null

# 7. Switch the page orientation to portrait for printing
$ws.PageSetup.Orientation = 1
